# Automated map update (2025-08-05 08:03:59)
# The INCO case with Caso=-522 (old row 43, "Uruguay 1090") was resolved/removed.
# Remove that entire row; the rows below it (old 44-46) shift up by one,
# which also shrinks the used range from P46 to P45.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A43").EntireRow.Delete()

Write-Host "Deleted row 43 (Caso -522 / Uruguay 1090); dimension shrinks from P46 to P45."
